$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 110, pushing the existing rows 110-172 down to 112-174.
$ws.Rows.Item(110).Insert()
$ws.Rows.Item(110).Insert()

# New row 110: "Primera" quality entry for the new date (44523).
$ws.Range("A110").Value = 11
$ws.Range("B110").Value = "Vega Monumental Concepción"
$ws.Range("C110").Value = "Bíobío"
$ws.Range("D110").Value = 44523
$ws.Range("E110").Value = 8
$ws.Range("F110").Value = 100112008
$ws.Range("G110").Value = "Coliflor"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 600
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 650
$ws.Range("N110").Value = "$/unidad"
$ws.Range("O110").Value = "Región Metropolitana"
$ws.Range("P110").Value = 650
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"

# New row 111: "Segunda" quality entry for the same new date (44523).
$ws.Range("A111").Value = 11
$ws.Range("B111").Value = "Vega Monumental Concepción"
$ws.Range("C111").Value = "Bíobío"
$ws.Range("D111").Value = 44523
$ws.Range("E111").Value = 8
$ws.Range("F111").Value = 100112008
$ws.Range("G111").Value = "Coliflor"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 500
$ws.Range("L111").Value = 500
$ws.Range("M111").Value = 500
$ws.Range("N111").Value = "$/unidad"
$ws.Range("O111").Value = "Región Metropolitana"
$ws.Range("P111").Value = 500
$ws.Range("Q111").Value = 1
$ws.Range("R111").Value = "Hortaliza"
